# Daily "remaining days" refresh for the water-delivery tracker sheet.
#
# Column D = total days (总天), E = days remaining (剩余), F = start date
# (开始时间, stored as an 8-digit yyyyMMdd integer).
#
# For every data row we recompute E as D minus the number of days elapsed
# since F, as of "today". If that countdown has reached zero (i.e. the
# delivery cycle finished), the row is refilled: F is reset to today and E
# goes back to the full D. Rows whose F value isn't a clean 8-digit date
# (data-entry errors) are left untouched, matching how the nightly updater
# skips rows it can't parse.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Today" for this run.
$todayDt = Get-Date -Year 2026 -Month 2 -Day 21 -Hour 0 -Minute 0 -Second 0
$todayOA = [math]::Floor($todayDt.ToOADate())
$todayNum = [int]($todayDt.ToString("yyyyMMdd"))

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {

    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value()
    $fVal = $fCell.Value()

    if ($null -eq $dVal -or $null -eq $fVal) {
        continue
    }

    $fStr = [string][int64]$fVal
    if ($fStr.Length -ne 8) {
        # Not a well-formed yyyyMMdd date (e.g. "202510929") - skip, same
        # as the source row being excluded from the automated refresh.
        continue
    }

    $y = [int]$fStr.Substring(0, 4)
    $mo = [int]$fStr.Substring(4, 2)
    $da = [int]$fStr.Substring(6, 2)

    $parsed = $true
    try {
        $fDt = Get-Date -Year $y -Month $mo -Day $da -Hour 0 -Minute 0 -Second 0
    } catch {
        $parsed = $false
    }
    if (-not $parsed) {
        continue
    }

    $fOA = [math]::Floor($fDt.ToOADate())
    $elapsed = $todayOA - $fOA
    $newE = [int]$dVal - $elapsed

    if ($newE -le 0) {
        # Cycle finished - restock: restart the countdown from today.
        $eCell.Value = [int]$dVal
        $fCell.Value = $todayNum
    } else {
        $eCell.Value = $newE
    }
}
